# Revert "Revert "Vers 1.1 BMC Gastroenterology""
#
# Adds a new "NOTES" worksheet (after the existing data sheet) that
# documents a handful of abbreviations used on the main sheet, and makes
# that new sheet the active/selected one.

$wb = $excel.ActiveWorkbook

# --- Add the new "NOTES" sheet after the last existing sheet -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$notes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$notes.Name = "NOTES"

# --- Populate the glossary table ------------------------------------------
$notes.Range("A1").Value = "ACM"
$notes.Range("B1").Value = "all cause mortality"

$notes.Range("A2").Value = "wm/bm/wf"
$notes.Range("B2").Value = "white or black; male or female"

$notes.Range("A3").Value = "mortality_xx_EA"
$notes.Range("B3").Value = "mortality from esophageal adenocarcinoma"

# --- Size the columns to fit their contents -------------------------------
$notes.Columns.Item(1).AutoFit()
$notes.Columns.Item(2).AutoFit()

# --- Leave the selection on the row below the table, as in the source ----
$notes.Range("A4").Select()
